$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()
$ws.Range("A1").Value = "TEST VALUE 4"
$ws.Protect()
Write-Host ($ws.Range("A1").Value2.ToString())
